# Automatische test-sync: 2025-06-29 15:22:50
# Appends the 19th test-mail entry to the "Logs" sheet and the matching
# category roll-up row to the "Dashboard" sheet, then grows the existing
# chart series ranges / conditional-formatting ranges to cover the new row.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Logs sheet: add row 34
# ---------------------------------------------------------------------
$logs = $wb.Worksheets.Item("Logs")

$logs.Range("A34").Value = "Zijn er vacatures?"
$logs.Range("B34").Value = "mailmind.test@zohomail.eu"
$logs.Range("C34").Value = "Testmail #19: Zijn er vacatures?"
$logs.Range("D34").Value = "Sollicitatie / Vacature"
$logs.Range("E34").Value = "Beste,`nDank u voor uw interesse in onze organisatie. Op dit moment hebben we geen openstaande vacatures, maar we moedigen u aan om regelmatig onze website te bezoeken voor eventuele toekomstige mogelijkheden. Mocht u nog vragen hebben of meer informatie wensen, aarzel dan niet om contact met ons op te nemen.`nMet vriendelijke groet,`n[Naam bedrijf]"
$logs.Range("F34").Value = "2025-06-29 15:22:21"
$logs.Range("G34").Value = "Ja"
$logs.Range("H34").Value = "Nee"
$logs.Range("I34").Value = "Ja"

# The existing conditional-formatting rules on D/G/H/I were scoped to
# row 33 as the last row; stretch each one down onto the new row 34 so
# the highlighted range keeps matching the data range (sheetData D2:D33
# -> D2:D34, etc.)
foreach ($col in @("D", "G", "H", "I")) {
    $oldRange = $logs.Range(($col + "2:" + $col + "33"))
    $newRange = $logs.Range(($col + "2:" + $col + "34"))
    $fcs = $oldRange.FormatConditions
    for ($i = 1; $i -le $fcs.Count; $i++) {
        $fcs.Item($i).ModifyAppliesToRange($newRange)
    }
}

# ---------------------------------------------------------------------
# 2. Dashboard sheet: add row 11 (new category tally)
# ---------------------------------------------------------------------
$dashboard = $wb.Worksheets.Item("Dashboard")

$dashboard.Range("A11").Value = "Sollicitatie / Vacature"
$dashboard.Range("B11").Value = 1

# ---------------------------------------------------------------------
# 3. Chart on the Dashboard sheet: grow the category/value series refs
#    from Dashboard!$A$2:$A$10 / $B$2:$B$10 to ...$11 so the new row
#    plots too.
# ---------------------------------------------------------------------
$chart = $dashboard.ChartObjects(1).Chart
$series = $chart.SeriesCollection(1)
$series.XValues = "='Dashboard'!`$A`$2:`$A`$11"
$series.Values = "='Dashboard'!`$B`$2:`$B`$11"
